$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Time Sheet")
$ws2 = $wb.Worksheets.Item("Protocol")

$hoursFormula = "=IFERROR(IF(COUNT(TimeSheet[[#This Row],[Time In]:[Time Out]])=4,(IF(TimeSheet[[#This Row],[Time Out]]<TimeSheet[[#This Row],[Time In]],1,0)+TimeSheet[[#This Row],[Time Out]])-TimeSheet[[#This Row],[Lunch End]]+TimeSheet[[#This Row],[Lunch Start]]-TimeSheet[[#This Row],[Time In]],IF(AND(LEN(TimeSheet[[#This Row],[Time In]])<>0,LEN(TimeSheet[[#This Row],[Time Out]])<>0),(IF(TimeSheet[[#This Row],[Time Out]]<TimeSheet[[#This Row],[Time In]],1,0)+TimeSheet[[#This Row],[Time Out]])-TimeSheet[[#This Row],[Time In]],0))*24,0)"

# ---------------------------------------------------------------------------
# Time Sheet: complete the previously half-filled row 21 (31.08.2018) and
# append new row 22 for 03.09.2018. Pull number formats from row 20, which
# already carries the correct per-column styling, so the new cells line up
# with the existing table styling instead of minting fresh style records.
# ---------------------------------------------------------------------------
$ws1.Range("C20:G20").Copy()
$ws1.Range("C21:G22").PasteSpecial(-4122)
$ws1.Range("A1").Select()

$ws1.Range("D21").Value = 0.510416666666667
$ws1.Range("E21").Value = 0.552083333333333
$ws1.Range("F21").Value = 0.708333333333333
$ws1.Range("G21").Formula = $hoursFormula

$ws1.Range("C22").Value = 0.385416666666667
$ws1.Range("D22").Value = 0.510416666666667
$ws1.Range("E22").Value = 0.552083333333333
$ws1.Range("F22").Value = 0.78125
$ws1.Range("G22").Formula = $hoursFormula

# B22 / H22 share column B / H's default ("General") style already, so a
# plain assignment keeps the implicit style - except dd.mm.yyyy-shaped text
# must be forced to text via a text-formatted staging cell + paste so Excel
# doesn't reinterpret it as a date serial.
$stage = $ws1.Range("K1")
$stage.NumberFormat = "@"
$stage.Value = "03.09.2018"
$stage.Copy()
$ws1.Range("B22").PasteSpecial(-4163)
$stage.Clear()

$ws1.Range("H22").Value = "Installing software for Kinect Camera"

# ---------------------------------------------------------------------------
# Protocol: new entries for 03.09.2018
# ---------------------------------------------------------------------------
$stage2 = $ws2.Range("K1")
$stage2.NumberFormat = "@"
$stage2.Value = "03.09.2018"
$stage2.Copy()
$ws2.Range("B49").PasteSpecial(-4163)
$stage2.ClearContents()

$ws2.Range("C49").Value = "Tried to install software for the Kinect camera at the workstation"
$ws2.Range("C50").Value = "There is no USB 3.0 at the PC so the Kinect can not be used " + [char]0x2192 + " setting everything up on a different computer"

# ---------------------------------------------------------------------------
# Final view state: Protocol tab active, selection on C51; Time Sheet's own
# last selection moves to H23 (below the newly added row).
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("H23").Select()

$ws2.Activate()
$ws2.Range("A25").Select()
$ws2.Range("C51").Select()

"done"
